# Weekly data refresh: insert two new price rows (Ciruela - Black Amber,
# week of 2022-01-28) at the top of the detail block (row 83), pushing the
# existing rows down by two. This mirrors the "Fruta / hortaliza, semanal"
# periodic append used by this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 83 (formatting is inherited from the row
# above, which already carries the date number-format on column D).
$ws.Rows("83:84").Insert()

# --- New row 83: Black Amber / Primera --------------------------------
$ws.Cells.Item(83, 1).Value = 8
$ws.Cells.Item(83, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(83, 3).Value = "Coquimbo"
$ws.Cells.Item(83, 4).Value = "1/28/2022"
$ws.Cells.Item(83, 5).Value = 4
$ws.Cells.Item(83, 6).Value = "Fruta"
$ws.Cells.Item(83, 7).Value = 100103
$ws.Cells.Item(83, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(83, 9).Value = 100103002
$ws.Cells.Item(83, 10).Value = "Ciruela"
$ws.Cells.Item(83, 11).Value = "Black Amber"
$ws.Cells.Item(83, 12).Value = "Primera"
$ws.Cells.Item(83, 13).Value = 20
$ws.Cells.Item(83, 14).Value = 255000
$ws.Cells.Item(83, 15).Value = 260000
$ws.Cells.Item(83, 16).Value = 257500
$ws.Cells.Item(83, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(83, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(83, 19).Value = 572
$ws.Cells.Item(83, 20).Value = 450

# --- New row 84: Black Amber / Segunda ---------------------------------
$ws.Cells.Item(84, 1).Value = 8
$ws.Cells.Item(84, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(84, 3).Value = "Coquimbo"
$ws.Cells.Item(84, 4).Value = "1/28/2022"
$ws.Cells.Item(84, 5).Value = 4
$ws.Cells.Item(84, 6).Value = "Fruta"
$ws.Cells.Item(84, 7).Value = 100103
$ws.Cells.Item(84, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(84, 9).Value = 100103002
$ws.Cells.Item(84, 10).Value = "Ciruela"
$ws.Cells.Item(84, 11).Value = "Black Amber"
$ws.Cells.Item(84, 12).Value = "Segunda"
$ws.Cells.Item(84, 13).Value = 20
$ws.Cells.Item(84, 14).Value = 205000
$ws.Cells.Item(84, 15).Value = 210000
$ws.Cells.Item(84, 16).Value = 207500
$ws.Cells.Item(84, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(84, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(84, 19).Value = 461
$ws.Cells.Item(84, 20).Value = 450
